$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.08785
$ws.Range("H2").Value = 0.26355
$ws.Range("M2").Value = 0.107627
$ws.Range("N2").Value = 0.322881
$ws.Range("O2").Value = 0.07374101335569673
$ws.Range("P2").Value = 0.07374101335569674
$ws.Range("Q2").Value = 0.009455031949999998
$ws.Range("R2").Value = 0.08509528755
$ws.Range("S2").Value = 0.07374101335569673
$ws.Range("T2").Value = 0.07374101335569674

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.08785
$ws.Range("H3").Value = 0.26355
$ws.Range("O3").Value = 0.7110803705584069
$ws.Range("P3").Value = 0.711080370558407
$ws.Range("Q3").Value = 0.09117433185
$ws.Range("R3").Value = 0.82056898665
$ws.Range("S3").Value = 0.7110803705584069
$ws.Range("T3").Value = 0.711080370558407

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.08785
$ws.Range("H4").Value = 0.26355
$ws.Range("M4").Value = 0.314059
$ws.Range("N4").Value = 0.942177
$ws.Range("O4").Value = 0.2151786160858963
$ws.Range("P4").Value = 0.2151786160858964
$ws.Range("Q4").Value = 0.02759008315
$ws.Range("R4").Value = 0.24831074835
$ws.Range("S4").Value = 0.2151786160858963
$ws.Range("T4").Value = 0.2151786160858964
